$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Refund"
